$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 into I1:J1 so the new header cells
# match the existing bold/centered/bordered header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row 2 is special: I2 = 6, J2 = 9 (H2 stays 4, unchanged)
$ws.Cells.Item(2, 9).Value = 6
$ws.Cells.Item(2, 10).Value = 9

# Rows 3-36: I = 1, J = value currently in H (H itself is unchanged)
for ($r = 3; $r -le 36; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
